$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/string updates (coin names, links, volume percentages)
$plainCells = @(
    "E2",
    "E3",
    "E4",
    "E5",
    "E6",
    "E7",
    "E8",
    "E9",
    "E10",
    "E11",
    "E12",
    "E13",
    "E14",
    "E16",
    "E17",
    "E18",
    "E19",
    "E20",
    "E21",
    "E23",
    "E24",
    "E25",
    "E26",
    "E27",
    "E28",
    "E29",
    "E30",
    "E31",
    "B32",
    "C32",
    "E32",
    "B33",
    "C33",
    "E33",
    "E34",
    "E35",
    "E36",
    "E37",
    "E38",
    "B39",
    "C39",
    "E39",
    "B40",
    "C40",
    "E40",
    "E41",
    "E42",
    "E43",
    "E44",
    "E45",
    "E46",
    "E47",
    "E48",
    "B49",
    "C49",
    "E49",
    "B50",
    "C50",
    "E50",
    "E51",
)
$plainValues = @(
    "  -2.99%  ",
    "  -3.36%  ",
    "  +0.14%  ",
    "  -3.06%  ",
    "  -6.66%  ",
    "  -3.39%  ",
    "  -0.10%  ",
    "  -5.99%  ",
    "  -9.16%  ",
    "  -3.89%  ",
    "  -10.52%  ",
    "  -6.99%  ",
    "  -3.24%  ",
    "  -2.99%  ",
    "  -3.01%  ",
    "  -1.26%  ",
    "  -9.36%  ",
    "  -10.95%  ",
    "  -8.25%  ",
    "  -9.05%  ",
    "  -5.53%  ",
    "  -3.18%  ",
    "  +0.05%  ",
    "  -11.58%  ",
    "  -13.39%  ",
    "  -10.62%  ",
    "  -10.59%  ",
    "  -13.94%  ",
    "ImmutableX",
    "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx",
    "  -9.66%  ",
    "NEARProtocol",
    "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near",
    "  -10.01%  ",
    "  +0.05%  ",
    "  -8.60%  ",
    "  -7.15%  ",
    "  -12.61%  ",
    "  -6.86%  ",
    "Stacks",
    "https://coinranking.com/coin/mMPrMcB7+stacks-stx",
    "  -8.18%  ",
    "USDe",
    "https://coinranking.com/coin/exbfr2U-0+usde-usde",
    "  -0.01%  ",
    "  -10.56%  ",
    "  +0.13%  ",
    "  -7.09%  ",
    "  -1.33%  ",
    "  -4.83%  ",
    "  -15.90%  ",
    "  -6.20%  ",
    "  -10.07%  ",
    "SuiNetwork",
    "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui",
    "  -5.19%  ",
    "InjectiveProtocol",
    "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj",
    "  -5.14%  ",
    "  -9.15%  ",
)
for ($i = 0; $i -lt $plainCells.Length; $i++) {
    $ws.Range($plainCells[$i]).Value = $plainValues[$i]
}

# Price column (D) updates - force text type to match original inlineStr formatting
$priceCells = @(
    "D2",
    "D3",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D15",
    "D16",
    "D17",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51",
)
$priceValues = @(
    "69.340.14",
    "3.683.66",
    "682.86",
    "159.95",
    "3.682.80",
    "0.999",
    "0.494",
    "0.146",
    "7.15",
    "0.436",
    "0.0000234",
    "4.303.67",
    "32.50",
    "3.670.76",
    "69.360.05",
    "15.88",
    "6.44",
    "473.50",
    "9.88",
    "0.651",
    "79.33",
    "3.824.28",
    "0.0000125",
    "10.94",
    "9.23",
    "2.69",
    "1.74",
    "2.03",
    "6.67",
    "26.69",
    "8.16",
    "6.07",
    "2.26",
    "1.00",
    "0.0906",
    "0.940",
    "165.02",
    "47.80",
    "2.72",
    "1.31",
    "0.000276",
    "1.10",
    "28.07",
    "7.89",
)
for ($i = 0; $i -lt $priceCells.Length; $i++) {
    $ws.Range($priceCells[$i]).NumberFormat = "@"
    $ws.Range($priceCells[$i]).Value = $priceValues[$i]
    $ws.Range($priceCells[$i]).Style = "Normal"
}
